# Scheduled runner update: refresh cached market-board price snapshots
# (currentAveragePrice / NQ / HQ / Leve price / profit columns H:N)
# across the per-job Leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 798.75
$ws.Range("I32").Value = 798.3333
$ws.Range("K32").Value = 798.3333
$ws.Range("M32").Value = -472.3333

$ws.Range("H99").Value = 1403.6364
$ws.Range("I99").Value = 415.8
$ws.Range("J99").Value = 2226.8333
$ws.Range("K99").Value = 1247.4
$ws.Range("L99").Value = 6680.499899999999
$ws.Range("M99").Value = 250.5999999999999
$ws.Range("N99").Value = -9676.499899999999

$ws.Range("H100").Value = 2388.2727
$ws.Range("I100").Value = 2136
$ws.Range("J100").Value = 2691
$ws.Range("K100").Value = 2136
$ws.Range("L100").Value = 2691
$ws.Range("M100").Value = -1595
$ws.Range("N100").Value = -3773

$ws.Range("H104").Value = 938.25
$ws.Range("I104").Value = 938.25
$ws.Range("K104").Value = 2814.75
$ws.Range("M104").Value = -1067.75

$ws.Range("H132").Value = 20490.076
$ws.Range("I132").Value = 19670.092
$ws.Range("K132").Value = 59010.276
$ws.Range("M132").Value = -56480.276

$ws.Range("H138").Value = 3167.3333
$ws.Range("I138").Value = 2070.1428
$ws.Range("J138").Value = 3619.1177
$ws.Range("K138").Value = 6210.428400000001
$ws.Range("L138").Value = 10857.3531
$ws.Range("M138").Value = -1070.428400000001
$ws.Range("N138").Value = -21137.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2768.889
$ws.Range("I45").Value = 1848.2
$ws.Range("K45").Value = 1848.2
$ws.Range("M45").Value = -1471.2

$ws.Range("H122").Value = 3036.8572
$ws.Range("I122").Value = 2998.8
$ws.Range("K122").Value = 8996.400000000001
$ws.Range("M122").Value = -6546.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 9695476
$ws.Range("I7").Value = 10454601
$ws.Range("K7").Value = 10454601
$ws.Range("M7").Value = -10454488

$ws.Range("H107").Value = 5114.5625
$ws.Range("I107").Value = 1166.625
$ws.Range("K107").Value = 1166.625
$ws.Range("M107").Value = 753.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 85.14286
$ws.Range("I7").Value = 66.083336
$ws.Range("K7").Value = 66.083336
$ws.Range("M7").Value = 46.916664

$ws.Range("H41").Value = 62191.25
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H50").Value = 5010
$ws.Range("I50").Value = 5010
$ws.Range("K50").Value = 5010
$ws.Range("M50").Value = -4385

$ws.Range("H88").Value = 9665.666999999999
$ws.Range("J88").Value = 9665.666999999999
$ws.Range("L88").Value = 9665.666999999999
$ws.Range("N88").Value = -10477.667

$ws.Range("H91").Value = 9665.666999999999
$ws.Range("J91").Value = 9665.666999999999
$ws.Range("L91").Value = 9665.666999999999
$ws.Range("N91").Value = -12473.667

$ws.Range("H107").Value = 1176.6666
$ws.Range("I107").Value = 967.2
$ws.Range("K107").Value = 967.2
$ws.Range("M107").Value = 952.8

$ws.Range("H132").Value = 4278.8
$ws.Range("I132").Value = 4278.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12836.4
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -10306.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 872222.25
$ws.Range("I4").Value = 872222.25
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2616666.75
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -2616554.75

$ws.Range("H37").Value = 119999
$ws.Range("J37").Value = 119998.75
$ws.Range("L37").Value = 359996.25
$ws.Range("N37").Value = -360220.25

$ws.Range("H107").Value = 489.64706
$ws.Range("I107").Value = 422
$ws.Range("J107").Value = 504.14285
$ws.Range("K107").Value = 1266
$ws.Range("L107").Value = 1512.42855
$ws.Range("M107").Value = 654
$ws.Range("N107").Value = -5352.428550000001

$ws.Range("H113").Value = 1777.1
$ws.Range("I113").Value = 1073.5
$ws.Range("J113").Value = 2246.1667
$ws.Range("K113").Value = 3220.5
$ws.Range("L113").Value = 6738.500100000001
$ws.Range("M113").Value = -1050.5
$ws.Range("N113").Value = -11078.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15364091
$ws.Range("I11").Value = 13125625
$ws.Range("J11").Value = 21333334
$ws.Range("K11").Value = 13125625
$ws.Range("L11").Value = 21333334
$ws.Range("M11").Value = -13125486
$ws.Range("N11").Value = -21333612

$ws.Range("H14").Value = 12361
$ws.Range("I14").Value = 300
$ws.Range("J14").Value = 15376.25
$ws.Range("K14").Value = 300
$ws.Range("L14").Value = 15376.25
$ws.Range("M14").Value = -132
$ws.Range("N14").Value = -15712.25

$ws.Range("H22").Value = 347.2
$ws.Range("I22").Value = 238.66667
$ws.Range("J22").Value = 510
$ws.Range("K22").Value = 238.66667
$ws.Range("L22").Value = 510
$ws.Range("M22").Value = 290.33333
$ws.Range("N22").Value = -1568

$ws.Range("H43").Value = 4627.9
$ws.Range("J43").Value = 5716.6665
$ws.Range("L43").Value = 5716.6665
$ws.Range("N43").Value = -6018.6665

$ws.Range("H122").Value = 2390.7778
$ws.Range("I122").Value = 2439.125
$ws.Range("K122").Value = 7317.375
$ws.Range("M122").Value = -4867.375

$ws.Range("H132").Value = 105031.2
$ws.Range("I132").Value = 203612.6
$ws.Range("J132").Value = 6449.8
$ws.Range("K132").Value = 610837.8
$ws.Range("L132").Value = 19349.4
$ws.Range("M132").Value = -608307.8
$ws.Range("N132").Value = -24409.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 15000
$ws.Range("I5").Value = 15000
$ws.Range("K5").Value = 15000
$ws.Range("M5").Value = -14887

$ws.Range("H22").Value = 995.375
$ws.Range("J22").Value = 1121.5
$ws.Range("L22").Value = 1121.5
$ws.Range("N22").Value = -1711.5

$ws.Range("H27").Value = 995.375
$ws.Range("J27").Value = 1121.5
$ws.Range("L27").Value = 1121.5
$ws.Range("N27").Value = -1335.5

$ws.Range("H46").Value = 4439
$ws.Range("J46").Value = 4708.0835
$ws.Range("L46").Value = 4708.0835
$ws.Range("N46").Value = -5084.0835

$ws.Range("H53").Value = 8186.75
$ws.Range("I53").Value = 8186.75
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 8186.75
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -7668.75

$ws.Range("H55").Value = 1244.4286
$ws.Range("I55").Value = 1386.8334
$ws.Range("J55").Value = 390
$ws.Range("K55").Value = 1386.8334
$ws.Range("L55").Value = 390
$ws.Range("M55").Value = -1213.8334
$ws.Range("N55").Value = -736

$ws.Range("H82").Value = 3425.7
$ws.Range("J82").Value = 4707.5
$ws.Range("L82").Value = 4707.5
$ws.Range("N82").Value = -5429.5

$ws.Range("H85").Value = 3425.7
$ws.Range("J85").Value = 4707.5
$ws.Range("L85").Value = 4707.5
$ws.Range("N85").Value = -7203.5

$ws.Range("H136").Value = 2362.5454
$ws.Range("I136").Value = 2362.5454
$ws.Range("K136").Value = 7087.6362
$ws.Range("M136").Value = -4537.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1206.1
$ws.Range("I96").Value = 1289
$ws.Range("J96").Value = 1012.6667
$ws.Range("K96").Value = 1289
$ws.Range("L96").Value = 1012.6667
$ws.Range("M96").Value = 84
$ws.Range("N96").Value = -3758.6667
